{"js": "// Center and bold the title paragraph (\"Java Key Words\"), and move the\n// \"_GoBack\" bookmark from the end of the document to the end of the\n// title paragraph (matches Word's behavior of tracking the last edit\n// location).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst title = paragraphs.items[0];\n\n// Center-align the paragraph and bold its text.\ntitle.alignment = Word.Alignment.centered;\ntitle.font.bold = true;\nawait context.sync();\n\n// Remove the existing \"_GoBack\" bookmark (currently at the end of the\n// document) and re-insert it at the end of the title paragraph.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst titleEnd = title.getRange(\"End\");\ntitleEnd.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Center and bold the title paragraph (\"Java Key Words\"), and move the\n# \"_GoBack\" bookmark from the end of the document to the end of the\n# title paragraph (matches Word's behavior of tracking the last edit\n# location).\n\n$d = $word.ActiveDocument\n\n$titlePara = $d.Paragraphs(1)\n\n# Center-align the paragraph and bold its text.\n$titlePara.Alignment = \"wdAlignParagraphCenter\"\n$titlePara.Range.Font.Bold = 1\n\n# Remove the existing \"_GoBack\" bookmark (currently at the end of the\n# document).\n$existing = $d.Bookmarks.Item(\"_GoBack\")\n$existing.Delete()\n\n# Re-insert \"_GoBack\" right after the title text, before the paragraph\n# mark. A bookmark collapsed exactly on a paragraph-end boundary is\n# ambiguous, so nudge it in using a temporary placeholder character that\n# is removed immediately after the bookmark is anchored.\n$titleContent = $titlePara.Range\n[void]$titleContent.MoveEnd(1, -1)\n$insertionPoint = $titleContent.Duplicate\n$insertionPoint.Collapse(0)\n$insertionPoint.InsertAfter(\"X\")\n\n$bookmarkRange = $d.Range($insertionPoint.Start, $insertionPoint.Start)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n\n$d.Range($insertionPoint.Start, $insertionPoint.Start + 1).Delete()\n"}
